# Apply the edit described by the commit:
#  - Remove the branding/logo pictures from several slide layouts and the
#    slide master (and the "TEKsystems ... RIGHTS RESERVED" textbox on the
#    master).
#  - Remove the stray "allegis logo" picture from slide 1 (the title slide).
#  - Delete the "Agenda" slide (slide id 267) entirely.

$p = $ppt.ActivePresentation

# --- Slide layouts: drop the small logo picture on each, keep the
#     full-bleed background picture where one is present. ---
$master = $p.Designs.Item(1).SlideMaster

# Layout 1 "Title Slide": single pic "Picture 6"
$cl = $master.CustomLayouts.Item(1)
$cl.Shapes.Item("Picture 6").Delete()

# Layout 2 "Title Slide 2": keep "Picture 4" (background), drop "Picture 8"
$cl = $master.CustomLayouts.Item(2)
$cl.Shapes.Item("Picture 8").Delete()

# Layout 3 "Section Header": single pic "Picture 5"
$cl = $master.CustomLayouts.Item(3)
$cl.Shapes.Item("Picture 5").Delete()

# Layout 4 "Section Header 2": keep "Picture 1" (background), drop "Picture 5"
$cl = $master.CustomLayouts.Item(4)
$cl.Shapes.Item("Picture 5").Delete()

# Layout 22 "Thank You": single pic "Picture 9"
$cl = $master.CustomLayouts.Item(22)
$cl.Shapes.Item("Picture 9").Delete()

# Layout 24 "2_Title and Content": single pic "Picture 4"
$cl = $master.CustomLayouts.Item(24)
$cl.Shapes.Item("Picture 4").Delete()

# --- Slide master: drop the logo picture + copyright textbox ---
$master.Shapes.Item("Picture 12").Delete()
$master.Shapes.Item("TextBox 13").Delete()

# --- Slide 1 (title slide): remove the allegis-logo picture ---
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item("Picture 2").Delete()

# --- Delete the "Agenda" slide (SlideID 267) ---
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    if ($s.SlideID -eq 267) {
        $s.Delete()
    }
}
